$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at the front (A and B); this shifts the
# existing columns (A..G) to (C..I), carrying their values/styles along.
$ws.Range("A:B").Insert()

# The new A1/B1 header cells should use the same (bold/bordered) header
# style as the rest of row 1 -- copy formats from the neighboring header
# cell, then set their text.
$ws.Range("C1").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)

$ws.Range("A1").Value = "data"
$ws.Range("B1").Value = "loja"

# Fill the new "data" (date) and "loja" (store) columns for every data row.
$ws.Range("A2").Value = "30/07/2024"
$ws.Range("B2").Value = "rbrondoneirelime"

$ws.Range("A3").Value = "30/07/2024"
$ws.Range("B3").Value = "rbrondoneirelime"

$ws.Range("A4").Value = "30/07/2024"
$ws.Range("B4").Value = "rbrondoneirelime"

$ws.Range("A5").Value = "30/07/2024"
$ws.Range("B5").Value = "rbrondoneirelime"

# Update the tracking_id query parameter inside the link column (now
# column I) for each row.
$ws.Range("I2").Value = "https://produto.mercadolivre.com.br/MLB-3687576338-fonte-carregador-jfa-storm-60-amperes-bivolt-com-voltimetro-_JM#position%3D1%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D1c6ba195-e19e-4b88-be84-7a934058ae19"
$ws.Range("I3").Value = "https://produto.mercadolivre.com.br/MLB-3334858103-fonte-carregador-jfa-60a-bivolt-storm-com-medidor-cca-_JM#position%3D2%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D1c6ba195-e19e-4b88-be84-7a934058ae19"
$ws.Range("I4").Value = "https://produto.mercadolivre.com.br/MLB-4376052480-filtro-antirruido-jfa-2020k-eletromagnetico-rca-para-som-_JM#position%3D3%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D1c6ba195-e19e-4b88-be84-7a934058ae19"
$ws.Range("I5").Value = "https://produto.mercadolivre.com.br/MLB-1800331683-filtro-anti-ruido-jfa-p-rca-cd-dvd-eletromagnetico-stereo-_JM#position%3D4%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D1c6ba195-e19e-4b88-be84-7a934058ae19"

Write-Host "Edit applied"
